# The document's two headers each carry the BTEC logo and its two
# footers each carry the Pearson logo, as inline pictures. The edit
# simply relabels the pictures' internal "name" (wp:docPr / pic:cNvPr
# name attribute, as seen e.g. in Word's Selection Pane) - the BTEC
# logo's name swaps from "image2.jpg" to "image1.jpg" and the Pearson
# logo's name swaps from "image1.png" to "image2.png". No visible
# content, size, or image data changes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers: BTEC logo, "image2.jpg" -> "image1.jpg" ---

try {
    $headerShape1 = $sec.Headers.Item(1).Range.InlineShapes.Item(1)
    $headerShape1.Name = "image1.jpg"
} catch {
    Write-Host "Header (primary) rename failed: $($_.Exception.Message)"
}

try {
    $headerShape2 = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
    $headerShape2.Name = "image1.jpg"
} catch {
    Write-Host "Header (first page) rename failed: $($_.Exception.Message)"
}

# --- Footers: Pearson logo, "image1.png" -> "image2.png" ---

try {
    $footerShape1 = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
    $footerShape1.Name = "image2.png"
} catch {
    Write-Host "Footer (primary) rename failed: $($_.Exception.Message)"
}

try {
    $footerShape2 = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
    $footerShape2.Name = "image2.png"
} catch {
    Write-Host "Footer (first page) rename failed: $($_.Exception.Message)"
}
